{"js": "// The document contains a paragraph built from three separate runs that,\n// concatenated, read: \"<id>\" + \"p167r_1\" + \"</id>\".\n// The edit merges that text into a single run (taking on the formatting\n// of the first run: Courier New / color 7f6000 / sz 18) so the paragraph\n// reads as one run: \"<id>p167r_1</id>\".\nconst body = context.document.body;\nconst results = body.search(\"<id>p167r_1</id>\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the '<id>p167r_1</id>' text to merge.\");\n}\n\n// insertText(\"Replace\") rewrites the matched range as a single run that\n// inherits the formatting of the range's leading run (the \"<id>\" run),\n// collapsing the three original runs into one.\nresults.items[0].insertText(\"<id>p167r_1</id>\", \"Replace\");\nawait context.sync();\n", "ps1": "# The document contains a paragraph built from three separate runs that,\n# concatenated, read: \"<id>\" + \"p167r_1\" + \"</id>\".\n# The edit merges that text into a single run (keeping the formatting of\n# the first run: Courier New / color 7f6000 / sz 18) so the paragraph\n# reads as one run: \"<id>p167r_1</id>\".\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$null = $rng.Find.Execute(\"<id>p167r_1</id>\")\n\nif (-not $rng.Find.Found) {\n    throw \"Could not find the '<id>p167r_1</id>' text to merge.\"\n}\n\n# Keep the leading \"<id>\" run untouched (it already carries the formatting\n# we want the merged run to end up with) and only touch the remainder, so\n# the re-inserted text lands inside that first run instead of minting a\n# brand-new, differently-seeded run.\n$rng.MoveStart(1, 4)  # wdCharacter = 1; \"<id>\" is 4 characters\n$rng.Delete()\n$rng.InsertAfter(\"p167r_1</id>\")\n"}
